$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-11 Thursday" "2025-12-12 Friday"

Replace-Text "12÷8=" "66÷4="
Replace-Text "56÷5=" "32÷7="
Replace-Text "41÷9=" "22÷8="
Replace-Text "62÷3=" "15÷7="
Replace-Text "81÷9=" "11÷9="
Replace-Text "43÷9=" "88÷5="
Replace-Text "44÷8=" "90÷8="
Replace-Text "38÷6=" "55÷7="
Replace-Text "25÷8=" "71÷2="
Replace-Text "21÷7=" "16÷6="
Replace-Text "35÷3=" "83÷5="
Replace-Text "20÷5=" "72÷4="
Replace-Text "84÷7=" "59÷7="
Replace-Text "33÷2=" "25÷5="
Replace-Text "90÷4=" "33÷3="
Replace-Text "22÷7=" "80÷2="
Replace-Text "29÷6=" "75÷4="
Replace-Text "14÷7=" "27÷7="
Replace-Text "43÷7=" "76÷4="
Replace-Text "63÷7=" "57÷5="
Replace-Text "65÷2=" "59÷6="
Replace-Text "33÷5=" "64÷2="
Replace-Text "13÷6=" "41÷4="
Replace-Text "50÷8=" "35÷5="
Replace-Text "52÷3=" "21÷7="

Write-Output "Done"
